# Boost.Alabaster deck — minor text/formatting edits on 3 slides.

$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# Slide 39: "A remedy for the review manager starvation." ->
#           "A Remedy for the Review Manager Starvation." (centered,
#           no bullet, re-cased/re-run text)
# -----------------------------------------------------------------
$s39 = $p.Slides.Item(39)
$sh39 = $s39.Shapes.Item(5)
$tr39 = $sh39.TextFrame.TextRange

$tr39.Text = "A Remedy for the Review Manager Starvation."
$tr39.Characters(1, 13).Text  = "A Remedy for "   # "A Remedy for "
$tr39.Characters(14, 11).Text = "the Review "     # "the Review "
$tr39.Characters(25, 1).Text  = "M"               # "M"
$tr39.Characters(26, 7).Text  = "anager "         # "anager "
$tr39.Characters(33, 1).Text  = "S"               # "S"
$tr39.Characters(34, 9).Text  = "tarvation"       # "tarvation"
$tr39.Characters(43, 1).Text  = "."               # "."

$tr39.ParagraphFormat.Alignment = 2          # ppAlignCenter
$tr39.ParagraphFormat.Bullet.Visible = $false

# -----------------------------------------------------------------
# Slide 40: "A desire to get the approval of experts" -> highlight
#           "approval" in blue (0070C0)
# -----------------------------------------------------------------
$s40 = $p.Slides.Item(40)
$sh40 = $s40.Shapes.Item(5)
$tr40 = $sh40.TextFrame.TextRange

$para40 = $tr40.Paragraphs(5, 1)
$approval = $tr40.Characters($para40.Start + 20, 8)
$approval.Font.Color.RGB = 0xC07000   # BGR for srgbClr 0070C0

# -----------------------------------------------------------------
# Slide 42: merge the split "oes all the work ... file a " / "final
#           report." runs back into a single run.
# -----------------------------------------------------------------
$s42 = $p.Slides.Item(42)
$sh42 = $s42.Shapes.Item(5)
$tr42 = $sh42.TextFrame.TextRange

$para42 = $tr42.Paragraphs(2, 1)
$tr42.Characters($para42.Start + 1, $para42.Length - 2).Text = "oes all the work that is necessary to check a library submission, organize the process, moderate and file a final report."
